$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 (the "Huge:" bullet, paragraph 9):
#   " ... so large that they are rather "  +  (italic) "moved"  +
#   " around than really carried - often only pushed ..."
#   becomes
#   " ... so large that they are "  +  (italic) "moved"  +
#   " around" + " rather" + " than carried - often only pushed ..."
# ---------------------------------------------------------------------------

$paraHuge = $d.Paragraphs.Item(9).Range
if ($paraHuge.Text -notlike "*Huge:*") {
    # fall back to a full-document search in case paragraph numbering differs
    $paraHuge = $d.Content
}

# 1a) Drop "really " before "carried" (isolated to the run after "moved", so
#     this cannot disturb the italic "moved" run or anything before it).
$r = $d.Paragraphs.Item(9).Range.Duplicate
$null = $r.Find.Execute(
    "than really carried", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "than carried", 2)

# 1b) Relocate the word "rather" (with its leading space) from right before
#     "moved" to right after "around". Cut/Paste (rather than a plain text
#     replace) keeps it as its own run on each end, matching how the
#     author's edit ends up splitting the sentence into separate runs.
$r = $d.Content
$null = $r.Find.Execute(
    " rather", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$r.Cut()

$r = $d.Content
$null = $r.Find.Execute(
    "around", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$dest = $d.Range($r.End, $r.End)
$dest.Paste()

# ---------------------------------------------------------------------------
# Change 2 (the "On Back:" bullet, paragraph 14):
#   "... such as a backpack, a " + [shield] + " or a large holstered weapon."
#   becomes
#   "... such as a backpack, a shield or a " + [large holstered] + " weapon."
#   (the run wrapped in <w:proofErr .../> must keep wrapping the swapped word)
# ---------------------------------------------------------------------------

$paraBack = $d.Paragraphs.Item(14).Range
if ($paraBack.Text -notlike "*On Back*") {
    $paraBack = $d.Content
}

# 2a) Trim the run that used to read " or a large holstered weapon." down to
#     just " weapon." -- do this first, while the phrase is still unique.
$r = $d.Paragraphs.Item(14).Range.Duplicate
$null = $r.Find.Execute(
    " or a large holstered weapon.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    " weapon.", 2)

# 2b) Swap the text inside the proofErr-wrapped run from "shield" to
#     "large holstered" (still the only occurrence of "shield" at this point).
$r = $d.Paragraphs.Item(14).Range.Duplicate
$null = $r.Find.Execute(
    "shield", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "large holstered", 2)

# 2c) Extend the run before it so it now ends with "... backpack, a shield or a ".
$r = $d.Paragraphs.Item(14).Range.Duplicate
$null = $r.Find.Execute(
    "such as a backpack, a ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "such as a backpack, a shield or a ", 2)
